$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 2 and row 3 (plus a few scattered cells further down) with the
# new NN-attribute sample data. Cells are written in the specific order
# below so that brand-new shared strings are interned in the same order
# as the target workbook: no, <800, 2013<model, >800, 2013>model, 60<=n, female.

# First occurrences of each new unique string, in target order
$ws.Range("F2").Value = "no"          # new string: no
$ws.Range("J2").Value = "<800"        # new string: <800
$ws.Range("K3").Value = "2013<model"  # new string: 2013<model
$ws.Range("J14").Value = ">800"       # new string: >800
$ws.Range("K15").Value = "2013>model" # new string: 2013>model
$ws.Range("A2").Value = "60<=n"       # new string: 60<=n
$ws.Range("B3").Value = "female"      # new string: female

# Remaining cells in row 2 (reusing existing shared strings / plain numbers)
# Note: B2/C2 already hold "male"/"none" with the correct styles in the
# original workbook, so they are intentionally left untouched here to avoid
# Excel resetting their cell style to the column default.
$ws.Range("D2").Value = 3000
$ws.Range("G2").Value = "no"
$ws.Range("H2").Value = "no"
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = "none"
$ws.Range("L2").Value = "none"

# Remaining cells in row 3
$ws.Range("A3").Value = "60<=n"
$ws.Range("C3").Value = "none"
$ws.Range("D3").Value = 3000
$ws.Range("E3").Value = "single"
$ws.Range("F3").Value = "no"
$ws.Range("G3").Value = "no"
$ws.Range("H3").Value = "no"
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = "<800"
$ws.Range("L3").Value = 1

# New cell in row 16
$ws.Range("L16").Value = 3

# Move the active selection to J1 (matches the saved view state)
[void]$ws.Range("J1").Select()
